$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 13 data (Day 12 entry)
$ws.Range("A13").Value = "Day 12"

# Copy the date cell's formatting from the row above so the new date cell
# reuses the existing date-number-format style instead of creating a new one,
# then set the serial date value (6/6/2025) without any time component.
$ws.Range("B12").Copy($ws.Range("B13"))
$ws.Range("B13").Value = 45814

$ws.Range("C13").Value = "Find the Index of the First Occurrence in a String"
$ws.Range("D13").Value = "Jump Game"
$ws.Range("E13").Value = "Length of Last Word"
$ws.Range("F13").Value = "Two Pointers, Greedy, String"
$ws.Range("G13").Value = "S"
$ws.Range("H13").Value = "YES"

# Widen column C slightly
$ws.Columns.Item(3).ColumnWidth = 38

# Update the selection to match the state after the edit
$ws.Range("B14").Select()
